$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Remove the "hasInterval" row (row 4) entirely; the blank row below shifts up to become row 4
$ws.Rows.Item(4).Delete()

# Update selection to reflect the new state
$ws.Range("A4:XFD4").Select()
